$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.123.84'
$ws.Range('E2').Value = '  +4.93%  '
$ws.Range('D3').Value = '2.741.91'
$ws.Range('E3').Value = '  +3.96%  '
$ws.Range('E4').Value = '  -0.27%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '580.14'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.42%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '158.02'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +10.01%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.625'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.53%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Value = '2.756.94'
$ws.Range('E9').Value = '  +3.66%  '
$ws.Range('E10').Value = '  +3.16%  '
$ws.Range('E11').Value = '  +3.75%  '
$ws.Range('E12').Value = '  +4.63%  '
$ws.Range('E13').Value = '  +1.25%  '
$ws.Range('D14').Value = '3.241.47'
$ws.Range('E14').Value = '  +3.79%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '27.33'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.66%  '
$ws.Range('D16').Value = '64.034.00'
$ws.Range('E16').Value = '  +4.72%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000156'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +7.80%  '
$ws.Range('D18').Value = '2.750.41'
$ws.Range('E18').Value = '  +3.18%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.13'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +5.18%  '
$ws.Range('E20').Value = '  +5.03%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '364.39'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.66%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.04'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.52%  '
$ws.Range('E23').Value = '  +4.28%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.998'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.17%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '66.96'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +4.66%  '
$ws.Range('E26').Value = '  +6.21%  '
$ws.Range('E27').Value = '  +1.82%  '
$ws.Range('E28').Value = '  +0.18%  '
$ws.Range('E29').Value = '  +13.66%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.02'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.09%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.27'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +7.04%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.26'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +15.40%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '174.10'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.77%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.999'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.09%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '20.66'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.57%  '
$ws.Range('E36').Value = '  +7.43%  '
$ws.Range('E37').Value = '  +10.04%  '
$ws.Range('E38').Value = '  +8.17%  '
$ws.Range('E39').Value = '  +12.66%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.29'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.36%  '
$ws.Range('B41').Value = 'Bittensor'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '338.77'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.18%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.18'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +16.38%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '39.55'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.14%  '
$ws.Range('B44').Value = 'InjectiveProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '22.56'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +7.71%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '21.91'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +7.26%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0606'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +5.61%  '
$ws.Range('E47').Value = '  +2.86%  '
$ws.Range('E48').Value = '  +3.98%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '137.92'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.99%  '
$ws.Range('E50').Value = '  +3.27%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.999'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.24%  '
